$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 no longer holds the old Hydrogen/Non-metallic-minerals value -> blank it out
$ws.Range("D3").ClearContents()

# Corrected code: C4 (Methanol/Chemicals) and C5 (Ammonia/Chemicals) values
$ws.Range("C4").Value = 207.0590250103807
$ws.Range("C5").Value = 6215.307221406939

# Row 7 label changes from "Other" to "Biogas", with an updated value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 4738.257339222111

# New row 8: re-introduce an "Other" row below Biogas, carrying the same
# label formatting as the other first-column entries (copy format from A7)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 3376.270858259668
